$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (legmaxROM header): columns B:E re-matched to the common max ROM set (15,16,15,16)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON passive force): recalculated meanEMG legmaxROM values for columns B:E
$ws.Range("B2").Value = 423.32962324447499
$ws.Range("C2").Value = 341.91203969887903
$ws.Range("D2").Value = 577.50735361712998
$ws.Range("E2").Value = 351.44220921357737

# Row 3 (STR passive force): recalculated meanEMG legmaxROM values for columns B:E
$ws.Range("B3").Value = 557.36065828433175
$ws.Range("C3").Value = 405.85417617272054
$ws.Range("D3").Value = 630.06752982188129
$ws.Range("E3").Value = 333.83415277465247

# Match the author's resulting selection (B1:E3) shown in the saved file
$ws.Range("B1:E3").Select() | Out-Null
